# The "Automation Testing" title slide carries a subtitle textbox with
# the speaker's name followed by a centered paragraph reading
# "   (Lead Product Efficiency Engineer)". That parenthetical job-title
# paragraph is being removed entirely (it disappears as a whole
# paragraph, not just blanked text).

$p = $ppt.ActivePresentation

$needle = "Lead Product Efficiency Engineer"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if (-not $shape.HasTextFrame) {
            continue
        }

        $textRange = $shape.TextFrame.TextRange

        if ($textRange.Text -notlike "*$needle*") {
            continue
        }

        # Walk paragraphs back-to-front so deleting one doesn't shift the
        # index of paragraphs we still need to inspect.
        for ($i = $textRange.Paragraphs().Count; $i -ge 1; $i--) {
            $para = $textRange.Paragraphs($i)
            if ($para.Text -like "*$needle*") {
                $para.Delete()
            }
        }
    }
}
